# Fruta / hortaliza, semanal
#
# Insert a new weekly record at row 485 in the "Espinaca" sheet, shifting
# the existing rows 485:509 down to 486:510 (dimension grows from R509 to
# R510). The new row carries a fresh weekly observation for the same
# market/product.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 485 and below down by one to make room for the new record.
$ws.Rows.Item(485).Insert()

# Populate the newly inserted row 485 with the new weekly observation.
$ws.Range("A485").Value = 8
$ws.Range("B485").Value = "Terminal La Palmera de La Serena"
$ws.Range("C485").Value = "Coquimbo"
$ws.Range("D485").Value = 45267
$ws.Range("E485").Value = 4
$ws.Range("F485").Value = 100112012
$ws.Range("G485").Value = "Espinaca"
$ws.Range("H485").Value = "Sin especificar"
$ws.Range("I485").Value = "Primera"
$ws.Range("J485").Value = 1300
$ws.Range("K485").Value = 450
$ws.Range("L485").Value = 500
$ws.Range("M485").Value = 475
$ws.Range("N485").Value = "`$/atado 300 a 500 gramos"
$ws.Range("O485").Value = "Provincia del Elquí"
$ws.Range("P485").Value = 950
$ws.Range("Q485").Value = 0.5
$ws.Range("R485").Value = "Hortaliza"
